$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix header row 2, column B: "unnamed: 1_level_1" -> "total"
# (the rest of row 2, columns C:I, already hold the correct labels)
$ws.Cells.Item(2, 2).Value = "total"

# Remove the two empty "category header" rows that no longer carry any
# data of their own ("situação do domicílio" and "grandes regiões e
# unidades da federação"). Deleting them shifts every row below up,
# which is exactly how the data realigns: "urbana"/"rondônia" (etc.)
# keep their own values, and the rows that used to hold only a label
# disappear.
$ws.Rows(5).Delete()
$ws.Rows(7).Delete()
